# @Cachelookup and webdriver fire event
#
# Refresh the "lastname" test-data values used by the Cachelookup test
# (King/Singh/Kumar -> King1/Singh1/Kumar1) and leave the grid selection
# on C4, matching where the webdriver-driven edit last landed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "King1"
$ws.Range("C3").Value = "Singh1"
$ws.Range("C4").Value = "Kumar1"

$ws.Range("C4").Select()
